# Refresh "想去人数" (want-to-go count, column F) figures across all sheets
# to match the latest upstream bilibili-huodong snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 525
$ws.Range("F6").Value = 945
$ws.Range("F9").Value = 1002
$ws.Range("F10").Value = 792
$ws.Range("F11").Value = 228
$ws.Range("F12").Value = 53
$ws.Range("F14").Value = 804
$ws.Range("F15").Value = 269
$ws.Range("F16").Value = 574
$ws.Range("F17").Value = 497
$ws.Range("F18").Value = 1321
$ws.Range("F20").Value = 576
$ws.Range("F21").Value = 1149
$ws.Range("F22").Value = 2840
$ws.Range("F23").Value = 1372
$ws.Range("F24").Value = 681
$ws.Range("F25").Value = 182
$ws.Range("F28").Value = 994
$ws.Range("F29").Value = 345
$ws.Range("F30").Value = 2532
$ws.Range("F31").Value = 459
$ws.Range("F32").Value = 521
$ws.Range("F33").Value = 1375

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 359
$ws.Range("F5").Value = 12
$ws.Range("F9").Value = 39

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 727

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 727
$ws.Range("F4").Value = 525
$ws.Range("F9").Value = 359
$ws.Range("F10").Value = 12
$ws.Range("F12").Value = 945
$ws.Range("F16").Value = 1002
$ws.Range("F17").Value = 792
$ws.Range("F18").Value = 228
$ws.Range("F20").Value = 53
$ws.Range("F21").Value = 39
$ws.Range("F26").Value = 804
$ws.Range("F27").Value = 269
$ws.Range("F28").Value = 574
$ws.Range("F29").Value = 497
$ws.Range("F30").Value = 1321
$ws.Range("F32").Value = 576
$ws.Range("F33").Value = 1149
$ws.Range("F34").Value = 2840
$ws.Range("F35").Value = 1372
$ws.Range("F36").Value = 681
$ws.Range("F37").Value = 182
$ws.Range("F42").Value = 994
$ws.Range("F43").Value = 345
$ws.Range("F44").Value = 2532
$ws.Range("F45").Value = 459
$ws.Range("F46").Value = 521
$ws.Range("F47").Value = 1375

Write-Output "Updated F-column attendance counts across sheets"
